$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3679.125
$ws.Range("I32").Value = 9800
$ws.Range("K32").Value = 9800
$ws.Range("M32").Value = -9474
$ws.Range("H43").Value = 3288.2856
$ws.Range("J43").Value = 3970.2
$ws.Range("L43").Value = 3970.2
$ws.Range("N43").Value = -4108.2
$ws.Range("H69").Value = 18922.154
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 18832.334
$ws.Range("K69").Value = 60000
$ws.Range("L69").Value = 56497.00199999999
$ws.Range("M69").Value = -59126
$ws.Range("N69").Value = -58245.00199999999
$ws.Range("H72").Value = 18922.154
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 18832.334
$ws.Range("K72").Value = 180000
$ws.Range("L72").Value = 169491.006
$ws.Range("M72").Value = -175632
$ws.Range("N72").Value = -178227.006
$ws.Range("H86").Value = 1404.9231
$ws.Range("I86").Value = 1496.4
$ws.Range("J86").Value = 1100
$ws.Range("K86").Value = 1496.4
$ws.Range("L86").Value = 1100
$ws.Range("M86").Value = -373.4000000000001
$ws.Range("N86").Value = -3346
$ws.Range("H89").Value = 1404.9231
$ws.Range("I89").Value = 1496.4
$ws.Range("J89").Value = 1100
$ws.Range("K89").Value = 7482
$ws.Range("L89").Value = 5500
$ws.Range("M89").Value = -1866
$ws.Range("N89").Value = -16732
$ws.Range("H131").Value = 4851.1
$ws.Range("I131").Value = 2387.2856
$ws.Range("K131").Value = 7161.8568
$ws.Range("M131").Value = -2121.8568
$ws.Range("H135").Value = 1564.3182
$ws.Range("I135").Value = 1482.579
$ws.Range("K135").Value = 13343.211
$ws.Range("M135").Value = -10808.211
$ws.Range("H137").Value = 2438.4
$ws.Range("I137").Value = 2267.5
$ws.Range("K137").Value = 6802.5
$ws.Range("M137").Value = -4252.5
$ws.Range("H138").Value = 3639.577
$ws.Range("J138").Value = 3225.3333
$ws.Range("L138").Value = 9675.999899999999
$ws.Range("N138").Value = -19955.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2448
$ws.Range("I45").Value = 1748.4445
$ws.Range("J45").Value = 3077.6
$ws.Range("K45").Value = 1748.4445
$ws.Range("L45").Value = 3077.6
$ws.Range("M45").Value = -1371.4445
$ws.Range("N45").Value = -3831.6
$ws.Range("H74").Value = 1982.6957
$ws.Range("I74").Value = 1384.6875
$ws.Range("K74").Value = 1384.6875
$ws.Range("M74").Value = -510.6875
$ws.Range("H77").Value = 1982.6957
$ws.Range("I77").Value = 1384.6875
$ws.Range("K77").Value = 6923.4375
$ws.Range("M77").Value = -2555.4375
$ws.Range("H97").Value = 2117.1155
$ws.Range("I97").Value = 351.45
$ws.Range("K97").Value = 351.45
$ws.Range("M97").Value = 144.55
$ws.Range("H110").Value = 3322.15
$ws.Range("I110").Value = 3247.111
$ws.Range("K110").Value = 3247.111
$ws.Range("M110").Value = -1202.111
$ws.Range("H122").Value = 3911.375
$ws.Range("I122").Value = 2907
$ws.Range("J122").Value = 4367.909
$ws.Range("K122").Value = 8721
$ws.Range("L122").Value = 13103.727
$ws.Range("M122").Value = -6271
$ws.Range("N122").Value = -18003.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1202.6316
$ws.Range("I99").Value = 1226.1428
$ws.Range("K99").Value = 1226.1428
$ws.Range("M99").Value = 271.8571999999999
$ws.Range("H107").Value = 4244.1665
$ws.Range("I107").Value = 3994.25
$ws.Range("J107").Value = 4744
$ws.Range("K107").Value = 3994.25
$ws.Range("L107").Value = 4744
$ws.Range("M107").Value = -2074.25
$ws.Range("N107").Value = -8584

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2733.5642
$ws.Range("I31").Value = 1758.6
$ws.Range("J31").Value = 3069.7585
$ws.Range("K31").Value = 1758.6
$ws.Range("L31").Value = 3069.7585
$ws.Range("M31").Value = -1463.6
$ws.Range("N31").Value = -3659.7585
$ws.Range("H34").Value = 2733.5642
$ws.Range("I34").Value = 1758.6
$ws.Range("J34").Value = 3069.7585
$ws.Range("K34").Value = 1758.6
$ws.Range("L34").Value = 3069.7585
$ws.Range("M34").Value = -1556.6
$ws.Range("N34").Value = -3473.7585

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 154
$ws.Range("I2").Value = 240.5
$ws.Range("J2").Value = 24.25
$ws.Range("K2").Value = 1443
$ws.Range("L2").Value = 145.5
$ws.Range("M2").Value = -1330
$ws.Range("N2").Value = -371.5
$ws.Range("H11").Value = 1829728.4
$ws.Range("I11").Value = 2195585.2
$ws.Range("J11").Value = 444.5
$ws.Range("K11").Value = 6586755.600000001
$ws.Range("L11").Value = 1333.5
$ws.Range("M11").Value = -6586615.600000001
$ws.Range("N11").Value = -1613.5
$ws.Range("H59").Value = 2668.3333
$ws.Range("I59").Value = 3002.5
$ws.Range("J59").Value = 2000
$ws.Range("K59").Value = 9007.5
$ws.Range("L59").Value = 6000
$ws.Range("M59").Value = -8467.5
$ws.Range("N59").Value = -7080
$ws.Range("H117").Value = 4414.375
$ws.Range("J117").Value = 4473.5713
$ws.Range("L117").Value = 13420.7139
$ws.Range("N117").Value = -20304.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1089.6316
$ws.Range("J97").Value = 2359.75
$ws.Range("L97").Value = 2359.75
$ws.Range("N97").Value = -3351.75
$ws.Range("H135").Value = 94749.5
$ws.Range("J135").Value = 94749.5
$ws.Range("L135").Value = 94749.5
$ws.Range("N135").Value = -104889.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 126870.25
$ws.Range("I7").Value = 168369.5
$ws.Range("K7").Value = 168369.5
$ws.Range("M7").Value = -168257.5
$ws.Range("H126").Value = 126870.25
$ws.Range("I126").Value = 168369.5
$ws.Range("K126").Value = 505108.5
$ws.Range("M126").Value = -502638.5
$ws.Range("H136").Value = 1179.8462
$ws.Range("I136").Value = 1153.1666
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3459.4998
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -909.4998000000001
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 440901.34
$ws.Range("I122").Value = 719201.3
$ws.Range("K122").Value = 2157603.9
$ws.Range("M122").Value = -2155153.9
$ws.Range("H126").Value = 1173.45
$ws.Range("I126").Value = 915
$ws.Range("K126").Value = 2745
$ws.Range("M126").Value = -275
$ws.Range("H136").Value = 9690.380999999999
$ws.Range("I136").Value = 12040.6455
$ws.Range("K136").Value = 36121.9365
$ws.Range("M136").Value = -33571.9365
